{"js": "// Briefing Final - Cliente.docx \u2014 field value updates (Office.js / Word JS API)\n//\n// The document is a flat list of paragraphs. Most paragraphs look like:\n//   <bold label run><br/><value run>\n// e.g. \"Endere\u00e7o completo:\" + line-break + \" N\u00e3o informado\"\n// Office.js renders the <w:br/> as a vertical-tab character (\"\\u000b\") inside\n// Paragraph.text, so a paragraph's text reads as \"Label:\\u000bValue\".\n//\n// Strategy: load every paragraph's text once, find the paragraph whose text\n// starts with the target label, then do a paragraph-scoped search() for the\n// OLD value text (so we never touch the many other paragraphs sharing the\n// same \"N\u00e3o informado\" placeholder) and insertText(..., \"Replace\") the hit.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// label -> [oldValue, newValue]\nconst fieldUpdates = [\n  [\"Endere\u00e7o completo:\", \" N\u00e3o informado\", \" asdasdas\"],\n  [\"Pr\u00eamios e reconhecimentos:\", \" N\u00e3o informado\", \" dasdasdas\"],\n  [\"Lista de produtos/servi\u00e7os:\", \" N\u00e3o informado\", \" dsadasdsa\"],\n  [\"Canais de compra:\", \" N\u00e3o informado\", \" asdasddas\"],\n  [\"Como clientes procuram:\", \" N\u00e3o informado\", \" dasdasdas\"],\n  [\"Concorrentes diretos:\", \" N\u00e3o informado\", \" dsadasdas\"],\n  [\"Melhor presen\u00e7a digital:\", \" dasdasdas\", \" N\u00e3o informado\"],\n  [\"Identidade visual:\", \" N\u00e3o informado\", \" Manual de marca completo\"],\n  [\"Tipos de fotos dispon\u00edveis:\", \" dasdsa\", \" N\u00e3o informado\"],\n  [\"Redes sociais:\", \" N\u00e3o informado\", \" dsadas\"],\n  [\"Plataformas presentes:\", \" dasdsadas\", \" N\u00e3o informado\"],\n  [\"Atributos do neg\u00f3cio:\", \" N\u00e3o informado\", \" dasdas\"],\n  [\"Caracter\u00edsticas do ambiente:\", \" dasdasdas\", \" N\u00e3o informado\"],\n  [\"Onde recebeu avalia\u00e7\u00f5es:\", \" N\u00e3o informado\", \" dewewawa\"],\n  [\"Meta de clientes mensais:\", \" N\u00e3o informado\", \" 1-10\"],\n  [\"Respons\u00e1vel pela gest\u00e3o:\", \" dsadasdas\", \" dsadsa\"],\n  [\"Tentativa anterior GMB:\", \" adsdasas\", \" dsadsa\"],\n  [\"Informa\u00e7\u00f5es a ocultar:\", \" dasdas\", \" dsadas\"],\n  [\"Restri\u00e7\u00f5es legais:\", \" N\u00e3o informado\", \" dasdas\"],\n  [\"Produtos restritos pelo Google:\", \" N\u00e3o informado\", \" dasdas\"],\n  [\"Maior expectativa:\", \" N\u00e3o informado\", \" dasdsa\"],\n];\n\n// Build label -> paragraph proxy map (first match wins; labels are unique).\nconst byLabel = {};\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  for (const [label] of fieldUpdates) {\n    if (!(label in byLabel) && text.indexOf(label) === 0) {\n      byLabel[label] = paragraphs.items[i];\n    }\n  }\n}\n\nfor (const [label, oldValue, newValue] of fieldUpdates) {\n  const para = byLabel[label];\n  if (!para) continue; // defensive: label not found, skip\n  const hits = para.search(oldValue, { matchCase: true });\n  hits.load(\"items\");\n  await context.sync();\n  if (hits.items.length > 0) {\n    hits.items[0].insertText(newValue, \"Replace\");\n  }\n}\n\n// Update the \"Data de envio\" timestamp line (2nd paragraph, plain text run \u2014\n// not a label/value pair, so handle it directly with a body-level search).\nconst dateHits = body.search(\"Data de envio: 24/06/2025, 17:44:47\", { matchCase: true });\ndateHits.load(\"items\");\nawait context.sync();\nif (dateHits.items.length > 0) {\n  dateHits.items[0].insertText(\"Data de envio: 24/06/2025, 18:17:28\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Briefing Final - Cliente.docx - field value updates (Word COM interop)\n#\n# The document is a flat list of paragraphs shaped like:\n#   <bold label run><line break><value run>\n# e.g. \"Endere\u00e7o completo:\" + line-break + \" N\u00e3o informado\"\n#\n# Strategy: walk $d.Paragraphs once, find the paragraph whose Range.Text\n# starts with the target label, then run Find/Replace scoped to THAT\n# paragraph's Range for the old value text (so the very common placeholder\n# \" N\u00e3o informado\" shared by many other paragraphs is never touched outside\n# the intended field).\n\n$d = $word.ActiveDocument\n\n$fieldUpdates = @(\n  @(\"Endere\u00e7o completo:\", \" N\u00e3o informado\", \" asdasdas\"),\n  @(\"Pr\u00eamios e reconhecimentos:\", \" N\u00e3o informado\", \" dasdasdas\"),\n  @(\"Lista de produtos/servi\u00e7os:\", \" N\u00e3o informado\", \" dsadasdsa\"),\n  @(\"Canais de compra:\", \" N\u00e3o informado\", \" asdasddas\"),\n  @(\"Como clientes procuram:\", \" N\u00e3o informado\", \" dasdasdas\"),\n  @(\"Concorrentes diretos:\", \" N\u00e3o informado\", \" dsadasdas\"),\n  @(\"Melhor presen\u00e7a digital:\", \" dasdasdas\", \" N\u00e3o informado\"),\n  @(\"Identidade visual:\", \" N\u00e3o informado\", \" Manual de marca completo\"),\n  @(\"Tipos de fotos dispon\u00edveis:\", \" dasdsa\", \" N\u00e3o informado\"),\n  @(\"Redes sociais:\", \" N\u00e3o informado\", \" dsadas\"),\n  @(\"Plataformas presentes:\", \" dasdsadas\", \" N\u00e3o informado\"),\n  @(\"Atributos do neg\u00f3cio:\", \" N\u00e3o informado\", \" dasdas\"),\n  @(\"Caracter\u00edsticas do ambiente:\", \" dasdasdas\", \" N\u00e3o informado\"),\n  @(\"Onde recebeu avalia\u00e7\u00f5es:\", \" N\u00e3o informado\", \" dewewawa\"),\n  @(\"Meta de clientes mensais:\", \" N\u00e3o informado\", \" 1-10\"),\n  @(\"Respons\u00e1vel pela gest\u00e3o:\", \" dsadasdas\", \" dsadsa\"),\n  @(\"Tentativa anterior GMB:\", \" adsdasas\", \" dsadsa\"),\n  @(\"Informa\u00e7\u00f5es a ocultar:\", \" dasdas\", \" dsadas\"),\n  @(\"Restri\u00e7\u00f5es legais:\", \" N\u00e3o informado\", \" dasdas\"),\n  @(\"Produtos restritos pelo Google:\", \" N\u00e3o informado\", \" dasdas\"),\n  @(\"Maior expectativa:\", \" N\u00e3o informado\", \" dasdsa\")\n)\n\n$count = $d.Paragraphs.Count\n\nforeach ($update in $fieldUpdates) {\n  $label = $update[0]\n  $oldValue = $update[1]\n  $newValue = $update[2]\n\n  for ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs($i)\n    $t = $p.Range.Text\n    if ($t.StartsWith($label)) {\n      $rng = $p.Range\n      $find = $rng.Find\n      $find.ClearFormatting()\n      $find.Text = $oldValue\n      $find.Replacement.ClearFormatting()\n      $find.Replacement.Text = $newValue\n      $find.Execute($null,$false,$false,$false,$false,$false,$true,1,$false,$null,2)\n      break\n    }\n  }\n}\n\n# Update the \"Data de envio\" timestamp line.\n$dateRange = $d.Content\n$dateFind = $dateRange.Find\n$dateFind.ClearFormatting()\n$dateFind.Text = \"Data de envio: 24/06/2025, 17:44:47\"\n$dateFind.Replacement.ClearFormatting()\n$dateFind.Replacement.Text = \"Data de envio: 24/06/2025, 18:17:28\"\n$dateFind.Execute($null,$false,$false,$false,$false,$false,$true,1,$false,$null,2)\n"}
